# billing keyword and page object update
# Adds a new "group" column (Y) with a "Dummy" value for every data row,
# and increases the height of rows 13-22 to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column Y (row 1)
$ws.Cells.Item(1, 25).Value = "group"

# "Dummy" placeholder value down column Y for every other row (2-45),
# including the blank separator row 23.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 25).Value = "Dummy"
}

# Rows 13-22 grow a bit taller to fit the new content.
for ($r = 13; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 24
}
